$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.67280729510103
$ws.Range("D2").Value = 8.566804068555699
$ws.Range("E2").Value = 14.73093161522375
$ws.Range("F2").Value = 37.16817263543215
$ws.Range("G2").Value = 44.26171340746591
$ws.Range("H2").Value = 17.82343068107678
$ws.Range("I2").Value = 27.02516683233621
$ws.Range("J2").Value = 10.86552837249869
$ws.Range("L2").Value = 12.92183536435359
$ws.Range("M2").Value = 18.02849632829411

$ws.Range("B3").Value = 17.26709732688305
$ws.Range("D3").Value = 8.465231213391602
$ws.Range("E3").Value = 14.49612687830674
$ws.Range("F3").Value = 37.27348896293255
$ws.Range("G3").Value = 44.10907488811738
$ws.Range("H3").Value = 17.87000092605768
$ws.Range("I3").Value = 27.24766214732932
$ws.Range("J3").Value = 10.77622715890346
$ws.Range("L3").Value = 12.71532515388489
$ws.Range("M3").Value = 17.79734882248679

$ws.Range("B4").Value = 17.01439893587058
$ws.Range("D4").Value = 8.40152559378604
$ws.Range("E4").Value = 14.34920236088944
$ws.Range("F4").Value = 37.35128231664274
$ws.Range("G4").Value = 44.03544366729268
$ws.Range("H4").Value = 17.90370112008706
$ws.Range("I4").Value = 27.39090549643653
$ws.Range("J4").Value = 10.7209460473496
$ws.Range("L4").Value = 12.58823988668032
$ws.Range("M4").Value = 17.65560099027812

$ws.Range("B5").Value = 16.91065004202817
$ws.Range("D5").Value = 8.375240514587116
$ws.Range("E5").Value = 14.28868306928892
$ws.Range("F5").Value = 37.38626375789543
$ws.Range("G5").Value = 44.01049746748826
$ws.Range("H5").Value = 17.91871189734099
$ws.Range("I5").Value = 27.45095136547477
$ws.Range("J5").Value = 10.69831722870684
$ws.Range("L5").Value = 12.53643628638564
$ws.Range("M5").Value = 17.59793751907775

$ws.Range("B6").Value = 16.89337990886347
$ws.Range("D6").Value = 8.370856660630462
$ws.Range("E6").Value = 14.27859621971044
$ws.Range("F6").Value = 37.39226985715975
$ws.Range("G6").Value = 44.00666076894236
$ws.Range("H6").Value = 17.92128140465937
$ws.Range("I6").Value = 27.46102315408627
$ws.Range("J6").Value = 10.69455395553126
$ws.Range("L6").Value = 12.52783508774188
$ws.Range("M6").Value = 17.58837019345007

$ws.Range("B7").Value = 17.01300269188429
$ws.Range("D7").Value = 8.401172401539245
$ws.Range("E7").Value = 14.34838873179544
$ws.Range("F7").Value = 37.351740837195
$ws.Range("G7").Value = 44.03508674699776
$ws.Range("H7").Value = 17.9038983960227
$ws.Range("I7").Value = 27.39170851434771
$ws.Range("J7").Value = 10.72064126183601
$ws.Range("L7").Value = 12.5875412316418
$ws.Range("M7").Value = 17.65482284301115

$ws.Range("B8").Value = 17.53374032343685
$ws.Range("D8").Value = 8.53206705604469
$ws.Range("E8").Value = 14.65056959480912
$ws.Range("F8").Value = 37.20174841720831
$ws.Range("G8").Value = 44.20492131869226
$ws.Range("H8").Value = 17.83842472800015
$ws.Range("I8").Value = 27.10051058968621
$ws.Range("J8").Value = 10.83483773909352
$ws.Range("L8").Value = 12.85072413952792
$ws.Range("M8").Value = 17.94879289674019

$ws.Range("B9").Value = 18.52098952962486
$ws.Range("D9").Value = 8.777598531975755
$ws.Range("E9").Value = 15.21935368572513
$ws.Range("F9").Value = 37.01267644864946
$ws.Range("G9").Value = 44.69660403036501
$ws.Range("H9").Value = 17.7508079006646
$ws.Range("I9").Value = 26.5818148201891
$ws.Range("J9").Value = 11.05469350867527
$ws.Range("L9").Value = 13.36209330000712
$ws.Range("M9").Value = 18.5242643127413

$ws.Range("B10").Value = 19.2188670303549
$ws.Range("D10").Value = 8.950514207585059
$ws.Range("E10").Value = 15.62025179224768
$ws.Range("F10").Value = 36.93894184721142
$ws.Range("G10").Value = 45.15291194557592
$ws.Range("H10").Value = 17.71163263210866
$ws.Range("I10").Value = 26.23226308940128
$ws.Range("J10").Value = 11.21305107262435
$ws.Range("L10").Value = 13.73171637095222
$ws.Range("M10").Value = 18.94326259167666

$ws.Range("B11").Value = 19.52910024305969
$ws.Range("D11").Value = 9.027417630849719
$ws.Range("E11").Value = 15.79845191475681
$ws.Range("F11").Value = 36.91975941555327
$ws.Range("G11").Value = 45.38060788478367
$ws.Range("H11").Value = 17.69934590503487
$ws.Range("I11").Value = 26.08001016398165
$ws.Range("J11").Value = 11.28426072508588
$ws.Range("L11").Value = 13.89790312216307
$ws.Range("M11").Value = 19.13239254272091

$ws.Range("B12").Value = 19.64543993190046
$ws.Range("D12").Value = 9.05627521223294
$ws.Range("E12").Value = 15.86529344538414
$ws.Range("F12").Value = 36.91457485165191
$ws.Range("G12").Value = 45.46966617221061
$ws.Range("H12").Value = 17.69549355220925
$ws.Range("I12").Value = 26.02332196582148
$ws.Range("J12").Value = 11.31109497182191
$ws.Range("L12").Value = 13.96050126104765
$ws.Range("M12").Value = 19.20374577172867

$ws.Range("B13").Value = 19.62043623871823
$ws.Range("D13").Value = 9.050072138576279
$ws.Range("E13").Value = 15.85092694850429
$ws.Range("F13").Value = 36.91559874695511
$ws.Range("G13").Value = 45.45036080896734
$ws.Range("H13").Value = 17.69628756181597
$ws.Range("I13").Value = 26.03548787556135
$ws.Range("J13").Value = 11.30532176923284
$ws.Range("L13").Value = 13.94703525010801
$ws.Range("M13").Value = 19.18839125996546

$ws.Range("B14").Value = 19.53869494168995
$ws.Range("D14").Value = 9.029797118043073
$ws.Range("E14").Value = 15.80396403164682
$ws.Range("F14").Value = 36.91929112751973
$ws.Range("G14").Value = 45.38787831184377
$ws.Range("H14").Value = 17.69901290490404
$ws.Range("I14").Value = 26.07532704881418
$ws.Range("J14").Value = 11.28647109564453
$ws.Range("L14").Value = 13.90306007458638
$ws.Range("M14").Value = 19.1382684546792

$ws.Range("B15").Value = 19.48847492522613
$ws.Range("D15").Value = 9.0173433591854
$ws.Range("E15").Value = 15.77511354327884
$ws.Range("F15").Value = 36.92182401595133
$ws.Range("G15").Value = 45.3499732876508
$ws.Range("H15").Value = 17.70078661363931
$ws.Range("I15").Value = 26.09985542384111
$ws.Range("J15").Value = 11.27490705424669
$ws.Range("L15").Value = 13.87607912224094
$ws.Range("M15").Value = 19.10753054602412

$ws.Range("B16").Value = 19.19843837929366
$ws.Range("D16").Value = 8.945452048711072
$ws.Range("E16").Value = 15.60851867013577
$ws.Range("F16").Value = 36.9404856999488
$ws.Range("G16").Value = 45.13843151985365
$ws.Range("H16").Value = 17.71254736189008
$ws.Range("I16").Value = 26.24234874150797
$ws.Range("J16").Value = 11.20837964256683
$ws.Range("L16").Value = 13.7208118214832
$ws.Range("M16").Value = 18.93086816816527

$ws.Range("B17").Value = 19.01858480415968
$ws.Range("D17").Value = 8.900890826519406
$ws.Range("E17").Value = 15.50522135994574
$ws.Range("F17").Value = 36.95562281003312
$ws.Range("G17").Value = 45.01377147435171
$ws.Range("H17").Value = 17.72118304891633
$ws.Range("I17").Value = 26.33149130707195
$ws.Range("J17").Value = 11.16734639538536
$ws.Range("L17").Value = 13.6250216280931
$ws.Range("M17").Value = 18.8220748347903

$ws.Range("B18").Value = 18.91446221769976
$ws.Range("D18").Value = 8.875095693013895
$ws.Range("E18").Value = 15.4454171583443
$ws.Range("F18").Value = 36.96568013561695
$ws.Range("G18").Value = 44.94396802979233
$ws.Range("H18").Value = 17.72667074093311
$ws.Range("I18").Value = 26.38340034765502
$ws.Range("J18").Value = 11.14366794323208
$ws.Range("L18").Value = 13.56974443319198
$ws.Range("M18").Value = 18.75936414579208

$ws.Range("B19").Value = 18.87909519811777
$ws.Range("D19").Value = 8.866333972789331
$ws.Range("E19").Value = 15.42510267599046
$ws.Range("F19").Value = 36.96931691357486
$ws.Range("G19").Value = 44.9206613672885
$ws.Range("H19").Value = 17.72861805900967
$ws.Range("I19").Value = 26.40108535724963
$ws.Range("J19").Value = 11.13563794056511
$ws.Range("L19").Value = 13.55099901830257
$ws.Range("M19").Value = 18.73810976405898

$ws.Range("B20").Value = 19.03780118113141
$ws.Range("D20").Value = 8.90565157438245
$ws.Range("E20").Value = 15.51625822267313
$ws.Range("F20").Value = 36.95387151962932
$ws.Range("G20").Value = 45.02684574103475
$ws.Range("H20").Value = 17.7202098433753
$ws.Range("I20").Value = 26.3219360830503
$ws.Range("J20").Value = 11.17172253511523
$ws.Range("L20").Value = 13.63523779488539
$ws.Range("M20").Value = 18.83367050673828

$ws.Range("B21").Value = 19.56273599876243
$ws.Range("D21").Value = 9.035759639222221
$ws.Range("E21").Value = 15.81777582350953
$ws.Range("F21").Value = 36.9181500449972
$ws.Range("G21").Value = 45.40615450359481
$ws.Range("H21").Value = 17.69819065016766
$ws.Range("I21").Value = 26.06359912522476
$ws.Range("J21").Value = 11.29201165554145
$ws.Range("L21").Value = 13.91598609331154
$ws.Range("M21").Value = 19.15299838488729

$ws.Range("B22").Value = 19.89913255731577
$ws.Range("D22").Value = 9.119248526729004
$ws.Range("E22").Value = 16.0110937376956
$ws.Range("F22").Value = 36.9069297794802
$ws.Range("G22").Value = 45.67054563013947
$ws.Range("H22").Value = 17.68846672621074
$ws.Range("I22").Value = 25.90039306914069
$ws.Range("J22").Value = 11.36985667816269
$ws.Range("L22").Value = 14.09750696577364
$ws.Range("M22").Value = 19.36012221696155

$ws.Range("B23").Value = 19.72023306701852
$ws.Range("D23").Value = 9.074833911618681
$ws.Range("E23").Value = 15.90827109705145
$ws.Range("F23").Value = 36.91180442821912
$ws.Range("G23").Value = 45.52794740724014
$ws.Range("H23").Value = 17.69322817386604
$ws.Range("I23").Value = 25.98698564706977
$ws.Range("J23").Value = 11.32838385161983
$ws.Range("L23").Value = 14.00082210856049
$ws.Range("M23").Value = 19.24973792866523

$ws.Range("B24").Value = 19.0291156995342
$ws.Range("D24").Value = 8.90349978856228
$ws.Range("E24").Value = 15.51126975358748
$ws.Range("F24").Value = 36.95465905903659
$ws.Range("G24").Value = 45.02092905060322
$ws.Range("H24").Value = 17.72064820106604
$ws.Range("I24").Value = 26.32625394419407
$ws.Range("J24").Value = 11.16974435519888
$ws.Range("L24").Value = 13.63061970397264
$ws.Range("M24").Value = 18.82842860996528

$ws.Range("B25").Value = 18.25827136015661
$ws.Range("D25").Value = 8.712440867232404
$ws.Range("E25").Value = 15.06829280944258
$ws.Range("F25").Value = 37.05245530252706
$ws.Range("G25").Value = 44.54674191780351
$ws.Range("H25").Value = 17.77010859624723
$ws.Range("I25").Value = 26.71657108092212
$ws.Range("J25").Value = 10.99572443256953
$ws.Range("L25").Value = 13.22459688284678
$ws.Range("M25").Value = 18.36902390008684
